# Release-Notes.xlsx update:
# A new folder "Azure_Well-Architected_Resiliency_Gaps_Remediation" was
# refreshed (its "Last Updated" timestamp moved forward), so it now sorts
# to the top of the "Folder Inventory" sheet. Its old row entry (row 47)
# is removed and a new row is inserted at row 2; every row that used to be
# between them (rows 2-46) shifts down by one row. Metadata / Summary
# sheets are updated accordingly.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Folder Inventory")

# Shift rows 2..46 down to rows 3..47 (work bottom-up so we never
# overwrite a row before it has been read).
for ($r = 46; $r -ge 2; $r--) {
    $dst = $r + 1
    $ws.Cells.Item($dst, 1).Value = $ws.Cells.Item($r, 1).Value2
    $ws.Cells.Item($dst, 2).Value = $ws.Cells.Item($r, 2).Value2
    $ws.Cells.Item($dst, 3).Value = $ws.Cells.Item($r, 3).Value2
    $ws.Cells.Item($dst, 4).Value = $ws.Cells.Item($r, 4).Value2
    $ws.Cells.Item($dst, 5).Value = $ws.Cells.Item($r, 5).Value2
}

# Write the refreshed entry into the now-empty row 2.
$ws.Cells.Item(2, 1).Value = "Azure_Well-Architected_Resiliency_Gaps_Remediation"
$ws.Cells.Item(2, 2).Value = "Azure_Well-Architected_Resiliency_Gaps_Remediation"
$ws.Cells.Item(2, 3).Value = "2025-06-13 17:35:45 +0530"
$ws.Cells.Item(2, 4).Value = 1
$ws.Cells.Item(2, 5).Value = "Root"

# --- Metadata sheet ---
$meta = $wb.Worksheets.Item("Metadata")
$meta.Cells.Item(3, 2).Value = "2025-06-13 12:06:07 UTC"

# "Workflow Run" (B5) is stored as text ("6"), not a number, so force the
# cell to text before assigning, then drop the temporary number format so
# no extra cell style is left behind.
$wfCell = $meta.Cells.Item(5, 2)
$wfCell.NumberFormat = "@"
$wfCell.Value = "6"
$wfCell.ClearFormats()

# --- Summary sheet ---
$summary = $wb.Worksheets.Item("Summary")
$summary.Cells.Item(5, 2).Value = "2025-06-13 17:35:45 +0530"
